$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.442379229604569
$ws.Range("C2").Value = 0.04231202262495515
$ws.Range("D2").Value = 0.1787246208189259
$ws.Range("E2").Value = 0.1640199045123722
$ws.Range("F2").Value = 1.578264622813109
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1910775826059847
$ws.Range("K2").Value = 0.4076749453714115
$ws.Range("O2").Value = 3.925773960152185
$ws.Range("B3").Value = 0.4028768011925195
$ws.Range("C3").Value = 0.03704680571253505
$ws.Range("D3").Value = 0.1726011840280535
$ws.Range("E3").Value = 0.1597661717657886
$ws.Range("F3").Value = 1.582086128929788
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1871676593352944
$ws.Range("K3").Value = 0.3655392386036453
$ws.Range("O3").Value = 3.949525278875086
$ws.Range("B4").Value = 0.3787075689788253
$ws.Range("C4").Value = 0.03379939988572289
$ws.Range("D4").Value = 0.1689111441244933
$ws.Range("E4").Value = 0.1572353702693761
$ws.Range("F4").Value = 1.585280313141325
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1848750716523568
$ws.Range("K4").Value = 0.3396974074448451
$ws.Range("O4").Value = 3.966251166974146
$ws.Range("B5").Value = 0.3688804300091135
$ws.Range("C5").Value = 0.03247247587924562
$ws.Range("D5").Value = 0.1674250754432052
$ws.Range("E5").Value = 0.1562244653391787
$ws.Range("F5").Value = 1.586795165354467
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1839680236687826
$ws.Range("K5").Value = 0.3291746831736475
$ws.Range("O5").Value = 3.973605674330926
$ws.Range("B6").Value = 0.3672499869534818
$ws.Range("C6").Value = 0.03225192722190684
$ws.Range("D6").Value = 0.1671793841980929
$ws.Range("E6").Value = 0.156057839890817
$ws.Range("F6").Value = 1.587059582515494
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1838190522281806
$ws.Range("K6").Value = 0.3274278948555889
$ws.Range("O6").Value = 3.974859408084583
$ws.Range("B7").Value = 0.378574946872277
$ws.Range("C7").Value = 0.03378151891641323
$ws.Range("D7").Value = 0.1688910308655096
$ws.Range("E7").Value = 0.1572216541202529
$ws.Range("F7").Value = 1.585299879721838
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1848627287371016
$ws.Range("K7").Value = 0.3395554610230818
$ws.Range("O7").Value = 3.966348172010854
$ws.Range("B8").Value = 0.4287413743881245
$ws.Range("C8").Value = 0.0404996382224283
$ws.Range("D8").Value = 0.176598839898972
$ws.Range("E8").Value = 0.1625364309417421
$ws.Range("F8").Value = 1.579406332969015
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1897070061629549
$ws.Range("K8").Value = 0.3931407266033773
$ws.Range("O8").Value = 3.933518635278148
$ws.Range("B9").Value = 0.527776739465537
$ws.Range("C9").Value = 0.05355572955383536
$ws.Range("D9").Value = 0.1922636791973389
$ws.Range("E9").Value = 0.1736000934153026
$ws.Range("F9").Value = 1.574576491666519
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.2000647204691148
$ws.Range("K9").Value = 0.4984368148981275
$ws.Range("O9").Value = 3.886149601908272
$ws.Range("B10").Value = 0.6009224178268653
$ws.Range("C10").Value = 0.06307327762640114
$ws.Range("D10").Value = 0.2041041074349579
$ws.Range("E10").Value = 0.1821188158426992
$ws.Range("F10").Value = 1.57513259589642
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2081990403758738
$ws.Range("K10").Value = 0.5759104939651252
$ws.Range("O10").Value = 3.861734003133762
$ws.Range("B11").Value = 0.6342785473484582
$ws.Range("C11").Value = 0.0673863040179441
$ws.Range("D11").Value = 0.2095618957747831
$ws.Range("E11").Value = 0.1860788460549543
$ws.Range("F11").Value = 1.576277746034719
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2120138081149179
$ws.Range("K11").Value = 0.6111762820214324
$ws.Range("O11").Value = 3.852885921632236
$ws.Range("B12").Value = 0.6469209721413449
$ws.Range("C12").Value = 0.06901709096939612
$ws.Range("D12").Value = 0.2116388137227858
$ws.Range("E12").Value = 0.187590571865762
$ws.Range("F12").Value = 1.576839710641721
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2134748211790196
$ws.Range("K12").Value = 0.6245333086606308
$ws.Range("O12").Value = 3.849860459283747
$ws.Range("B13").Value = 0.6441977116143107
$ws.Range("C13").Value = 0.06866598264291213
$ws.Range("D13").Value = 0.2111910616671366
$ws.Range("E13").Value = 0.1872644551427243
$ws.Range("F13").Value = 1.576712973964177
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2131594347291923
$ws.Range("K13").Value = 0.6216565264168139
$ws.Range("O13").Value = 3.850497582169567
$ws.Range("B14").Value = 0.6353184264416427
$ws.Range("C14").Value = 0.06752051966317651
$ws.Range("D14").Value = 0.2097325618455272
$ws.Range("E14").Value = 0.1862029734567372
$ws.Range("F14").Value = 1.576321407300867
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.2121336770538278
$ws.Range("K14").Value = 0.6122751232828421
$ws.Range("O14").Value = 3.852630497445745
$ws.Range("B15").Value = 0.6298810470244405
$ws.Range("C15").Value = 0.06681856739201919
$ws.Range("D15").Value = 0.2088405099514148
$ws.Range("E15").Value = 0.1855543658939212
$ws.Range("F15").Value = 1.576098273409158
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2115075119383505
$ws.Range("K15").Value = 0.6065290683346518
$ws.Range("O15").Value = 3.853979318458101
$ws.Range("B16").Value = 0.5987441126297881
$ws.Range("C16").Value = 0.06279107106995241
$ws.Range("D16").Value = 0.2037488576571604
$ws.Range("E16").Value = 0.1818617211947284
$ws.Range("F16").Value = 1.575075712846612
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2079520379316335
$ws.Range("K16").Value = 0.5736061943566142
$ws.Range("O16").Value = 3.862357723876926
$ws.Range("B17").Value = 0.5796631481686347
$ws.Range("C17").Value = 0.06031603404329644
$ws.Range("D17").Value = 0.2006435314558814
$ws.Range("E17").Value = 0.1796180931908751
$ws.Range("F17").Value = 1.574676932284333
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2058001690943314
$ws.Range("K17").Value = 0.5534144331505217
$ws.Range("O17").Value = 3.868076335990111
$ws.Range("B18").Value = 0.5686960153505538
$ws.Range("C18").Value = 0.0588909049383517
$ws.Range("D18").Value = 0.198864168283464
$ws.Range("E18").Value = 0.1783356043504725
$ws.Range("F18").Value = 1.574531538433803
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.2045732419241375
$ws.Range("K18").Value = 0.5418028320913777
$ws.Range("O18").Value = 3.87157811079058
$ws.Range("B19").Value = 0.5649840785095819
$ws.Range("C19").Value = 0.05840811644661414
$ws.Range("D19").Value = 0.1982628672247273
$ws.Range("E19").Value = 0.1779027488643266
$ws.Range("F19").Value = 1.574496732381206
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.2041596756995432
$ws.Range("K19").Value = 0.5378717342827883
$ws.Range("O19").Value = 3.872800252813221
$ws.Range("B20").Value = 0.5816935532510854
$ws.Range("C20").Value = 0.06057966731411568
$ws.Range("D20").Value = 0.2009734021542187
$ws.Range("E20").Value = 0.1798561049879837
$ws.Range("F20").Value = 1.57471069194068
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2060281245396567
$ws.Range("K20").Value = 0.5555636609009866
$ws.Range("O20").Value = 3.867445576395994
$ws.Range("B21").Value = 0.637926188849093
$ws.Range("C21").Value = 0.06785703738947291
$ws.Range("D21").Value = 0.2101606834280147
$ws.Range("E21").Value = 0.1865144270462693
$ws.Range("F21").Value = 1.576432937108876
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2124345207685678
$ws.Range("K21").Value = 0.6150306020812195
$ws.Range("O21").Value = 3.851995182507636
$ws.Range("B22").Value = 0.6747423798062471
$ws.Range("C22").Value = 0.07259882822219765
$ws.Range("D22").Value = 0.2162243408204176
$ws.Range("E22").Value = 0.1909368123763997
$ws.Range("F22").Value = 1.578306489725179
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2167173100589821
$ws.Range("K22").Value = 0.6539107793159644
$ws.Range("O22").Value = 3.843792616769093
$ws.Range("B23").Value = 0.6550871370456548
$ws.Range("C23").Value = 0.0700693889727404
$ws.Range("D23").Value = 0.2129826695103532
$ws.Range("E23").Value = 0.1885700412287505
$ws.Range("F23").Value = 1.577238091639529
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2144227391767686
$ws.Range("K23").Value = 0.6331585183457094
$ws.Range("O23").Value = 3.847996968988951
$ws.Range("B24").Value = 0.5807755984494065
$ws.Range("C24").Value = 0.06046048555123207
$ws.Range("D24").Value = 0.2008242491321255
$ws.Range("E24").Value = 0.1797484767929944
$ws.Range("F24").Value = 1.574695167948462
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2059250340771257
$ws.Range("K24").Value = 0.5545920046106403
$ws.Range("O24").Value = 3.867730076011213
$ws.Range("B25").Value = 0.5009162959379978
$ws.Range("C25").Value = 0.05003664357617765
$ws.Range("D25").Value = 0.1879674467498802
$ws.Range("E25").Value = 0.1705384921364654
$ws.Range("F25").Value = 1.575162526404611
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1971706747313675
$ws.Range("K25").Value = 0.4699303097603718
$ws.Range("O25").Value = 3.897141435657403
